$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 14.768619
$ws.Range("H2").Value = 44.305857
$ws.Range("I2").Value = 0.9736679609684162
$ws.Range("J2").Value = 0.9736679609684162
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 112.513392
$ws.Range("N2").Value = 337.540176
$ws.Range("O2").Value = 0.3275312977368564
$ws.Range("P2").Value = 0.3275312977368564
$ws.Range("Q2").Value = 1661.667418845648
$ws.Range("R2").Value = 14955.00676961083
$ws.Range("S2").Value = 0.3189067308207842
$ws.Range("T2").Value = 0.3189067308207842

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 14.768619
$ws.Range("H3").Value = 44.305857
$ws.Range("I3").Value = 0.9736679609684162
$ws.Range("J3").Value = 0.9736679609684162
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 106.314466
$ws.Range("N3").Value = 318.943398
$ws.Range("O3").Value = 0.3094859589441663
$ws.Range("P3").Value = 0.3094859589441664
$ws.Range("Q3").Value = 1570.117842542454
$ws.Range("R3").Value = 14131.06058288209
$ws.Range("S3").Value = 0.3013365625935214
$ws.Range("T3").Value = 0.3013365625935214

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 14.768619
$ws.Range("H4").Value = 44.305857
$ws.Range("I4").Value = 0.9736679609684162
$ws.Range("J4").Value = 0.9736679609684162
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 124.6916553333333
$ws.Range("N4").Value = 374.074966
$ws.Range("O4").Value = 0.3629827433189773
$ws.Range("P4").Value = 0.3629827433189773
$ws.Range("Q4").Value = 1841.523550097318
$ws.Range("R4").Value = 16573.71195087586
$ws.Range("S4").Value = 0.3534246675541106
$ws.Range("T4").Value = 0.3534246675541107

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.399405
$ws.Range("H5").Value = 1.198215
$ws.Range("I5").Value = 0.02633203903158381
$ws.Range("J5").Value = 0.02633203903158381
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 112.513392
$ws.Range("N5").Value = 337.540176
$ws.Range("O5").Value = 0.3275312977368564
$ws.Range("P5").Value = 0.3275312977368564
$ws.Range("Q5").Value = 44.93841133176
$ws.Range("R5").Value = 404.44570198584
$ws.Range("S5").Value = 0.008624566916072199
$ws.Range("T5").Value = 0.008624566916072199

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.399405
$ws.Range("H6").Value = 1.198215
$ws.Range("I6").Value = 0.02633203903158381
$ws.Range("J6").Value = 0.02633203903158381
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 106.314466
$ws.Range("N6").Value = 318.943398
$ws.Range("O6").Value = 0.3094859589441663
$ws.Range("P6").Value = 0.3094859589441664
$ws.Range("Q6").Value = 42.46252929273
$ws.Range("R6").Value = 382.16276363457
$ws.Range("S6").Value = 0.008149396350644931
$ws.Range("T6").Value = 0.008149396350644933

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.399405
$ws.Range("H7").Value = 1.198215
$ws.Range("I7").Value = 0.02633203903158381
$ws.Range("J7").Value = 0.02633203903158381
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 124.6916553333333
$ws.Range("N7").Value = 374.074966
$ws.Range("O7").Value = 0.3629827433189773
$ws.Range("P7").Value = 0.3629827433189773
$ws.Range("Q7").Value = 49.80247059841
$ws.Range("R7").Value = 448.2222353856901
$ws.Range("S7").Value = 0.009558075764866678
$ws.Range("T7").Value = 0.00955807576486668
